$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.01253208636536152; C = 0.04103571897497393; D = 0.7210945179870265;  E = 13.86384647080068; F = 0; G = 14.63850879412805 }
    3 = @{ B = 0.6545652718822623;  C = 1.626987699542094;  D = 0.7210945179870265;  E = 0.5333859586016987; F = 0; G = 3.536033448013082 }
    4 = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 3.223369029078222;   E = 0.5333859586016987; F = 0; G = 8.656069925401464 }
    5 = @{ B = 0.01253208636536152; C = 0.3048912486333797; D = 18.71679738969934;   E = 13.86384647080068; F = 0; G = 32.89806719549876 }
    6 = @{ B = 1.445647641019636;   C = 1.626987699542094;  D = 3.223369029078222;   E = 0.5333859586016987; F = 0; G = 6.82939032824165 }
    7 = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 0.1496068669990043;  E = 0.5333859586016987; F = 0; G = 5.582307763322248 }
    8 = @{ B = 0.1169995834814548;  C = 0.04103571897497393; D = 0.7210945179870265;  E = 0.5333859586016987; F = 0; G = 1.412515779045154 }
    9 = @{ B = 0.1169995834814548;  C = 0.3048912486333797; D = 0.1496068669990043;  E = 0.5333859586016987; F = 1; G = 1.104883657715537 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
}

$wb.Save()
